$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the header block (course / author / project title / date) plus a
#    blank line before the existing first paragraph ("Company Application
#    Management...").
# ---------------------------------------------------------------------------
$origFirstCount = $d.Paragraphs.Count
$headRange = $d.Range(0, 0)
$headRange.InsertBefore("CS5200`rRandy Lirano`rProject 1: Application Manager`r07/04/2021`r")

$addedLines = $d.Paragraphs.Count - $origFirstCount
$firstExisting = $d.Paragraphs($addedLines + 1)
$firstExisting.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 2) Mark the runs that hold the two inline pictures as NoProof (adds
#    <w:rPr><w:noProof/></w:rPr> to those runs), matching what Word stamps
#    on a freshly (re)inserted picture.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shapeRange = $d.InlineShapes.Item($i).Range
    $shapeRange.NoProofing = 1
}

# ---------------------------------------------------------------------------
# 3) Fix the mis-split "applicant_id" column name: the stray leading "a"
#    currently lives at the end of the previous run (", a") instead of at
#    the start of "pplicant_id".
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("pplicant_id)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$brokenStart = $findRange.Start

$commaRun = $d.Range($brokenStart - 3, $brokenStart)
$commaRun.Text = ", "

$wordRun = $d.Range($brokenStart - 1, $brokenStart - 1 + 11)
$wordRun.Text = "applicant_id"

Write-Output "done"
